$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows above the current row 11 ("feet"/"ft") to make room for
# the new "verify"/"return"/"build" key-value pairs.
$ws.Range("A11:D13").EntireRow.Insert()

$ws.Range("A11").Value = "verify"
$ws.Range("B11").Value = "VERIFY"

$ws.Range("A12").Value = "return"
$ws.Range("B12").Value = "RETURN"

$ws.Range("A13").Value = "build"
$ws.Range("B13").Value = "BUILD"

# Match the wrap-text style used by sibling rows in column B.
$ws.Range("B11:B13").WrapText = $true

# Restore default top-left view and move selection to B12, matching the
# saved workbook state.
$ws.Range("B12").Select()
$excel.ActiveWindow.ScrollRow = 1
